# Applies text corrections / updates to the UC005 test-suite workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Version bump 0.1 -> 1.2.5
$ws.Range("D2").Value = "1.2.5"

# Precondition text correction (accent + trailing period) - repeated in every test case block
$precondition = "O usuário devidamente autenticado e na tela de listagem de empenhos."
$ws.Range("B8").Value  = $precondition
$ws.Range("B16").Value = $precondition
$ws.Range("B24").Value = $precondition
$ws.Range("B31").Value = $precondition
$ws.Range("B38").Value = $precondition
$ws.Range("B45").Value = $precondition
$ws.Range("B53").Value = $precondition

# "filtra" -> "Filtra" (capitalised) + trailing period - repeated step text
$filtra = "Chefe/Beneficiário Filtra a listagem por registros cujos beneficiários não possuem número do credor."
$ws.Range("B10").Value = $filtra
$ws.Range("B47").Value = $filtra
$ws.Range("B55").Value = $filtra

# TC2 expected result: added "de todos os servidores," + accent fixes
$ws.Range("D18").Value = "SYSTEM Exibe a lista de solicitações aguardando serem empenhadas, de todos os servidores, ordenado pelo número da diária em ordem crescente."

# TC3 / TC4 step content swapped: "realizar o empenho" step now comes first (TC3),
# "atribuir/desatribuir" step now comes second (TC4)
$ws.Range("B26").Value = "Chefe/Beneficiário Clica para realizar o empenho de uma diária."
$ws.Range("D26").Value = "SYSTEM Apresenta a tela de Registrar Empenho."

$ws.Range("B33").Value = "Chefe/Beneficiário Clica para atribuir/desatribuir o registro a si mesmo."
$ws.Range("D33").Value = "SYSTEM Atualiza a lista de registros de solicitações, onde o nome deverá constar o nome do usuário logado (que se atribuiu como responsável pelo empenho), no campo de atribuição (no caso de desatribuição, o nome deverá ser removido)."

# TC5 expected result: trailing period added
$ws.Range("D40").Value = "SYSTEM Recupera e exibe todos os detalhes (dados) da solicitação para o usuário; e Apresenta a tela de Detalhar Diárias."
